# "added defend base state"
#
# The DefendBase rule sentence:
#   "... if the base health is under 20% and is being attacked by an enemy
#   tank, "
# becomes:
#   "... if the base health is under 20% and " / "or" / " being attacked by
#   an enemy tank"
# (i.e. "is" -> "or", and the trailing ", " is dropped) and the "_GoBack"
# bookmark is moved from the empty paragraph that follows onto the end of
# this (now shorter) sentence.

$d = $word.ActiveDocument

$full = $d.Content.Text
$oldSentence = "is being attacked by an enemy tank, "
$idx = $full.IndexOf($oldSentence)

if ($idx -ge 0) {
    # Boundary right before "is" -- keeps the untouched text ahead of it
    # (ending in "and ") as its own run instead of being re-merged with the
    # freshly typed "or".
    $boundaryBeforeIs = $d.Range($idx, $idx)
    $d.Bookmarks.Add("ZZ_TMP_SPLIT_1", $boundaryBeforeIs)

    # Boundary right after "is" -- likewise keeps the text that follows
    # (" being attacked by an enemy tank") as its own run.
    $boundaryAfterIs = $d.Range($idx + 2, $idx + 2)
    $d.Bookmarks.Add("ZZ_TMP_SPLIT_2", $boundaryAfterIs)

    # "is" -> "or"
    $isRange = $d.Range($idx, $idx + 2)
    $isRange.Text = "or"

    # Drop the trailing ", " after "...enemy tank"
    $full2 = $d.Content.Text
    $tailTarget = " being attacked by an enemy tank, "
    $tidx = $full2.IndexOf($tailTarget)
    if ($tidx -ge 0) {
        $tailRange = $d.Range($tidx + $tailTarget.Length - 2, $tidx + $tailTarget.Length)
        $tailRange.Delete()
    } else {
        Write-Host "WARNING: trailing ', ' not found for removal"
    }

    # The temporary bookmarks have done their job of keeping the runs apart;
    # remove them again now.
    $d.Bookmarks.Item("ZZ_TMP_SPLIT_1").Delete()
    $d.Bookmarks.Item("ZZ_TMP_SPLIT_2").Delete()
} else {
    Write-Host "WARNING: target sentence not found, no replacement made"
}

# Move the "_GoBack" bookmark so it sits right after the sentence we just
# edited, instead of in the (still) empty paragraph that follows it.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$full3 = $d.Content.Text
$anchor = "and or being attacked by an enemy tank"
$aidx = $full3.IndexOf($anchor)
if ($aidx -ge 0) {
    $pos = $aidx + $anchor.Length
    $bookmarkRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)
} else {
    Write-Host "WARNING: could not locate insertion point for _GoBack bookmark"
}
